# GodotCommonResource.xlsx test-fixture update:
#  - B29: 170 -> 150
#  - B30: 110 -> 100
#  - new row 31: key "speedBiliBili" / value 200
#  - leave selection on C32 (matches the saved cursor position in the diff)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B29").Value = 150
$ws.Range("B30").Value = 100

$ws.Range("A31").Value = "speedBiliBili"
$ws.Range("B31").Value = 200

$ws.Range("C32").Select()
